# Update "想去人数" (want-to-go count) values for two events that are
# duplicated across the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览": row 3 = 南宁·三月三漫次元国风动漫节, row 4 = 南宁·2024三月三国潮动漫节（良牙春典）
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 361
$wsExhibit.Range("F4").Value = 2960

# Sheet "全部类型": row 5 = 南宁·三月三漫次元国风动漫节, row 6 = 南宁·2024三月三国潮动漫节（良牙春典）
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 361
$wsAll.Range("F6").Value = 2960
